$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This workbook is a weekly price log. Two new weekly observation rows were
# added to the "Mandarina" sheet:
#   - a new row inserted at row 216 (pushing the old 216..238 down to 217..239)
#   - a new row inserted at row 235 (pushing what is now 235..239 down to 236..240)
# All other rows keep their original values; they just move down by one or
# two positions. We reproduce that with native row Insert() (shift down),
# then populate the two freshly inserted rows with the new data.
# ---------------------------------------------------------------------------

# 1) Insert a new blank row at 216 (existing rows 216-238 shift to 217-239)
$ws.Rows.Item(216).Insert()

# 2) Insert a second new blank row at 235 (existing rows 235-239 shift to 236-240)
$ws.Rows.Item(235).Insert()

# 3) Fill in the new row 216
$ws.Cells.Item(216,1).Value = 4
$ws.Cells.Item(216,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(216,3).Value = "Los Lagos"
$ws.Cells.Item(216,4).Value = 44748
$ws.Cells.Item(216,5).Value = 10
$ws.Cells.Item(216,6).Value = "Fruta"
$ws.Cells.Item(216,7).Value = 100102
$ws.Cells.Item(216,8).Value = "Cítricos"
$ws.Cells.Item(216,9).Value = 100102004
$ws.Cells.Item(216,10).Value = "Mandarina"
$ws.Cells.Item(216,11).Value = "Clemenuless"
$ws.Cells.Item(216,12).Value = "Primera"
$ws.Cells.Item(216,13).Value = 400
$ws.Cells.Item(216,14).Value = 9000
$ws.Cells.Item(216,15).Value = 10000
$ws.Cells.Item(216,16).Value = 9500
$ws.Cells.Item(216,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(216,18).Value = "Provincia de Limarí"
$ws.Cells.Item(216,19).Value = 950
$ws.Cells.Item(216,20).Value = 10

# 4) Fill in the new row 235
$ws.Cells.Item(235,1).Value = 4
$ws.Cells.Item(235,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(235,3).Value = "Los Lagos"
$ws.Cells.Item(235,4).Value = 44747
$ws.Cells.Item(235,5).Value = 10
$ws.Cells.Item(235,6).Value = "Fruta"
$ws.Cells.Item(235,7).Value = 100102
$ws.Cells.Item(235,8).Value = "Cítricos"
$ws.Cells.Item(235,9).Value = 100102004
$ws.Cells.Item(235,10).Value = "Mandarina"
$ws.Cells.Item(235,11).Value = "Clemenuless"
$ws.Cells.Item(235,12).Value = "Primera"
$ws.Cells.Item(235,13).Value = 800
$ws.Cells.Item(235,14).Value = 9000
$ws.Cells.Item(235,15).Value = 9000
$ws.Cells.Item(235,16).Value = 9000
$ws.Cells.Item(235,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(235,18).Value = "Provincia de Limarí"
$ws.Cells.Item(235,19).Value = 900
$ws.Cells.Item(235,20).Value = 10

# 5) Make sure the date columns keep the original date/time number format
#    (row Insert() normally inherits formatting from the row above already,
#    but set it explicitly for the two new date cells to be safe).
$ws.Cells.Item(216,4).NumberFormat = $ws.Cells.Item(217,4).NumberFormat
$ws.Cells.Item(235,4).NumberFormat = $ws.Cells.Item(236,4).NumberFormat
